# Updates the cryptos price/volume table (Sheet1) with the latest scraped
# values, matching the "Updated cryptos list ... with GitHub Actions" commit.
# Cells whose new text looks like a plain number (single decimal point, e.g.
# "210.93") are forced to the Text number format first so Excel keeps them
# as strings instead of converting them to numeric values, matching the
# original inline-string cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.588.10'
$ws.Range('E2').Value = '  -2.14%  '
$ws.Range('D3').Value = '1.583.62'
$ws.Range('E3').Value = '  -2.82%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.93'
$ws.Range('E5').Value = '  -2.41%  '
$ws.Range('E6').Value = '  -2.09%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.248'
$ws.Range('E8').Value = '  -2.53%  '
$ws.Range('E9').Value = '  -0.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.48'
$ws.Range('E10').Value = '  -3.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0833'
$ws.Range('E11').Value = '  -1.76%  '
$ws.Range('D12').Value = '1.806.04'
$ws.Range('E12').Value = '  -2.84%  '
$ws.Range('D13').Value = '1.581.81'
$ws.Range('E13').Value = '  -3.01%  '
$ws.Range('E14').Value = '  -1.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.527'
$ws.Range('E15').Value = '  -2.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.62'
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').Value = '26.619.88'
$ws.Range('E17').Value = '  -1.99%  '
$ws.Range('E18').Value = '  -0.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '208.11'
$ws.Range('E19').Value = '  -3.26%  '
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.75'
$ws.Range('E21').Value = '  -2.12%  '
$ws.Range('E22').Value = '  -3.05%  '
$ws.Range('E23').Value = '  -3.82%  '
$ws.Range('E24').Value = '  -2.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.13'
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.42'
$ws.Range('E26').Value = '  +2.03%  '
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('E28').Value = '  -4.27%  '
$ws.Range('E29').Value = '  -1.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0502'
$ws.Range('E30').Value = '  -0.41%  '
$ws.Range('E31').Value = '  -2.33%  '
$ws.Range('E33').Value = '  +23.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.93'
$ws.Range('E34').Value = '  -2.76%  '
$ws.Range('D35').Value = '1.317.60'
$ws.Range('E35').Value = '  +0.41%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.51'
$ws.Range('E36').Value = '  -3.56%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.42'
$ws.Range('E37').Value = '  -5.32%  '
$ws.Range('E38').Value = '  -1.12%  '
$ws.Range('E39').Value = '  -3.16%  '
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.785'
$ws.Range('E41').Value = '  -2.18%  '
$ws.Range('E42').Value = '  +2.33%  '
$ws.Range('E43').Value = '  -3.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.40'
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('D45').Value = '1.719.41'
$ws.Range('E45').Value = '  -2.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.11'
$ws.Range('E46').Value = '  -1.75%  '
$ws.Range('E47').Value = '  +1.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.831'
$ws.Range('E48').Value = '  +3.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0993'
$ws.Range('E49').Value = '  +4.47%  '
$ws.Range('E50').Value = '  -1.74%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.49'
$ws.Range('E51').Value = '  -0.39%  '
